$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update indicator values (B2:B13) ---
$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsMetrics.Range("B2").Value = 350661.08000000007
$wsMetrics.Range("B3").Value = 308941.79000000004
$wsMetrics.Range("B4").Value = 108037.02
$wsMetrics.Range("B5").Value = 14300
$wsMetrics.Range("B6").Value = 5146906.830000001
$wsMetrics.Range("B7").Value = 4351018.4700000007
$wsMetrics.Range("B8").Value = 1514996.8500000003
$wsMetrics.Range("B9").Value = 200507
$wsMetrics.Range("B10").Value = 33612287.820000015
$wsMetrics.Range("B11").Value = 31626293.630000006
$wsMetrics.Range("B12").Value = 11796718.889999999
$wsMetrics.Range("B13").Value = 1298137

# Move the Metrics selection (it is no longer the active/displayed sheet)
$wsMetrics.Range("E17").Select() | Out-Null

# --- today sheet: clear the "yesterday" helper values (B3:B6) ---
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Range("B3").ClearContents()
$wsToday.Range("B4").ClearContents()
$wsToday.Range("B5").ClearContents()
$wsToday.Range("B6").ClearContents()

# Make "today" the active sheet/tab with the new selection
$wsToday.Activate() | Out-Null
$wsToday.Range("F8").Select() | Out-Null
